$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clan games")

# --- Row 17: update E and J values ---
$ws.Range("E17").Value = 45900
$ws.Range("J17").Value = 9300

# --- Row 38: update E and J values ---
$ws.Range("E38").Value = 76150
$ws.Range("J38").Value = 4450

# --- Rows 40-47: overwrite with the new (shifted) participant data ---
$ws.Range("A40").Value = "Dasters79"
$ws.Range("B40").Value = "21/02/2026"
$ws.Range("E40").Value = 27155
$ws.Range("J40").Value = 4250

$ws.Range("A41").Value = "dibba10"
$ws.Range("B41").Value = "21/02/2026"
$ws.Range("E41").Value = 13360
$ws.Range("J41").Value = 900

$ws.Range("A42").Value = "cucco"
$ws.Range("B42").Value = "21/02/2026"
$ws.Range("E42").Value = 29400
$ws.Range("J42").Value = 8050

$ws.Range("A43").Value = "fede61mito"
$ws.Range("B43").Value = "21/02/2026"
$ws.Range("E43").Value = 400
$ws.Range("J43").Value = 0

$ws.Range("A44").Value = "MIRIAM MIRIAM"
$ws.Range("B44").Value = "21/02/2026"
$ws.Range("E44").Value = 23750
$ws.Range("J44").Value = 4300

$ws.Range("A45").Value = "PESCARA MANZIA"
$ws.Range("B45").Value = "23/02/2026"
$ws.Range("E45").Value = 11775
$ws.Range("J45").Value = 0

$ws.Range("A46").Value = "VERRETHERULER"
$ws.Range("B46").Value = "23/02/2026"
$ws.Range("E46").Value = 21145
$ws.Range("J46").Value = 0

$ws.Range("A47").Value = "GGfresco_08"
$ws.Range("B47").Value = "24/02/2026"
$ws.Range("E47").Value = 2150
$ws.Range("J47").Value = 0

# --- Rows 48-50: these participants no longer exist, remove entirely ---
$ws.Range("A48:AA48").Clear()
$ws.Range("A49:AA49").Clear()
$ws.Range("A50:AA50").Clear()

# --- Row 51: clear data/content but keep the row's own formatting (border/height) ---
$ws.Range("A51:AA51").Clear()
